# "error when loading map" fix
# Corrects the TQ (column F) values for several Spartan Heavy Infantry (HI)
# units on the "Spartan" sheet, which were all mistakenly set to 8.
# Also removes stray leftover helper-formula rows (16-21) that only ever
# evaluated to 0 and were no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spartan")

# Avoid automatic recalculation side effects (e.g. pivot tables silently
# refreshing) while we make the edits; we'll recalculate only what is
# needed afterwards.
$excel.Calculation = -4135  # xlCalculationManual

# TQ column (F) corrections: GHI1..GHI7 (rows 2-8) go from 8 to 6,
# SHI1..SHI2 (rows 9-10) go from 8 to 7.
$ws.Range("F2:F8").Value2 = 6
$ws.Range("F9:F10").Value2 = 7

# Refresh the dependent SUM formula's cached value without forcing a full
# workbook recalculation (which would also needlessly refresh pivot tables).
$ws.Range("F33").Formula = "=SUM(F2:F30)"

# Remove the now-unused helper rows (their L/M formulas only ever
# evaluated to 0 and had no other data).
$ws.Range("L16:M21").ClearContents()

# Update the active selection to reflect where the edit was made.
[void]$ws.Range("F8").Select()
